$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty rows 10-13 with new journal entries.
# Shared-string pool order follows first-use, so write B13 (whose text
# became the first newly-added shared string) before the others.
$ws.Range("B13").Value = "Discussion planification et organisation de groupe"
$ws.Range("B10").Value = "Lecture du tutoriel JavaFX sur openclassroom"
$ws.Range("B11").Value = "Visionnage de vidéos sur JavaFX et FXML sur la chaîne youtube thenewboston"
$ws.Range("B12").Value = "Création du projet à l'aide d'une borderPane ainsi que création de la MenuBar"

$ws.Range("A10").Value = "03/07/2018"
$ws.Range("C10").Value = 3

$ws.Range("A11").Value = "03/10/2018"
$ws.Range("C11").Value = 2

$ws.Range("A12").Value = "03/14/2018"
$ws.Range("C12").Value = 4

$ws.Range("A13").Value = "03/16/2018"
$ws.Range("C13").Value = 2

# The longer wrapped activity text in rows 11-12 makes Excel grow those
# rows to fit the content.
$ws.Rows.Item(11).RowHeight = 28.5
$ws.Rows.Item(12).RowHeight = 33

# Update selection to match the author's saved cursor position.
$ws.Range("C11").Select()
